$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.958.07'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.555.88'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.24'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.12'
$ws.Range('E8').Value = '  +4.00%  '
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0858'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').Value = '1.778.64'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '1.556.64'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').Value = '26.969.48'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.72'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '218.13'
$ws.Range('E18').Value = '  +2.30%  '
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.07'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.95'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').Value = '1.421.96'
$ws.Range('E33').Value = '  +4.74%  '
$ws.Range('E34').Value = '  +5.03%  '
$ws.Range('E35').Value = '  +4.10%  '
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.523'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.812'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.76'
$ws.Range('E41').Value = '  +2.71%  '
$ws.Range('E42').Value = '  -0.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.32'
$ws.Range('E43').Value = '  +4.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.985'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.39'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('E46').Value = '  +1.46%  '
$ws.Range('D47').Value = '1.691.89'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.09'
$ws.Range('E48').Value = '  +2.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0521'
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('D50').Value = '0.0₆01000'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0956'
$ws.Range('E51').Value = '  +1.02%  '
